$wb = $excel.ActiveWorkbook

# Sheet ALC row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 664.36365  # H12: was 729
$ws.Cells.Item(12, 9).Value = 590.4286  # I12: was 627.1667
$ws.Cells.Item(12, 10).Value = 793.75  # J12: was 932.6667
$ws.Cells.Item(12, 11).Value = 590.4286  # K12: was 627.1667
$ws.Cells.Item(12, 12).Value = 793.75  # L12: was 932.6667
$ws.Cells.Item(12, 13).Value = -420.4286  # M12: was -457.1667
$ws.Cells.Item(12, 14).Value = -1133.75  # N12: was -1272.6667

# Sheet ALC row 39
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(39, 8).Value = 1123.7  # H39: was 1318.1428
$ws.Cells.Item(39, 9).Value = 1176.8572  # I39: was 1371.3334
$ws.Cells.Item(39, 10).Value = 999.6667  # J39: was 999
$ws.Cells.Item(39, 11).Value = 3530.5716  # K39: was 4114.0002
$ws.Cells.Item(39, 12).Value = 2999.0001  # L39: was 2997
$ws.Cells.Item(39, 13).Value = -3234.5716  # M39: was -3818.0002
$ws.Cells.Item(39, 14).Value = -3591.0001  # N39: was -3589

# Sheet ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1799.2222  # H137: was 1768.3
$ws.Cells.Item(137, 9).Value = 1497.25  # I137: was 1495.8
$ws.Cells.Item(137, 11).Value = 4491.75  # K137: was 4487.4
$ws.Cells.Item(137, 13).Value = -1941.75  # M137: was -1937.4

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3132.7  # H32: was 3156.3667
$ws.Cells.Item(32, 9).Value = 2014.1305  # I32: was 2045
$ws.Cells.Item(32, 11).Value = 2014.1305  # K32: was 2045
$ws.Cells.Item(32, 13).Value = -1727.1305  # M32: was -1758

# Sheet ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1728.5  # H45: was 1711.2858
$ws.Cells.Item(45, 9).Value = 1643.25  # I45: was 1636.2
$ws.Cells.Item(45, 11).Value = 1643.25  # K45: was 1636.2
$ws.Cells.Item(45, 13).Value = -1266.25  # M45: was -1259.2

# Sheet ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 3120.611  # H74: was 3298.2942
$ws.Cells.Item(74, 9).Value = 3104.6667  # I74: was 3204.6667
$ws.Cells.Item(74, 10).Value = 3200.3333  # J74: was 4000.5
$ws.Cells.Item(74, 11).Value = 3104.6667  # K74: was 3204.6667
$ws.Cells.Item(74, 12).Value = 3200.3333  # L74: was 4000.5
$ws.Cells.Item(74, 13).Value = -2230.6667  # M74: was -2330.6667
$ws.Cells.Item(74, 14).Value = -4948.3333  # N74: was -5748.5

# Sheet ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 3120.611  # H77: was 3298.2942
$ws.Cells.Item(77, 9).Value = 3104.6667  # I77: was 3204.6667
$ws.Cells.Item(77, 10).Value = 3200.3333  # J77: was 4000.5
$ws.Cells.Item(77, 11).Value = 15523.3335  # K77: was 16023.3335
$ws.Cells.Item(77, 12).Value = 16001.6665  # L77: was 20002.5
$ws.Cells.Item(77, 13).Value = -11155.3335  # M77: was -11655.3335
$ws.Cells.Item(77, 14).Value = -24737.6665  # N77: was -28738.5

# Sheet ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1784.9524  # H132: was 1824.9445
$ws.Cells.Item(132, 9).Value = 1006.1429  # I132: was 1065.4615
$ws.Cells.Item(132, 10).Value = 3342.5715  # J132: was 3799.6
$ws.Cells.Item(132, 11).Value = 3018.4287  # K132: was 3196.3845
$ws.Cells.Item(132, 12).Value = 10027.7145  # L132: was 11398.8
$ws.Cells.Item(132, 13).Value = -488.4287000000004  # M132: was -666.3844999999997
$ws.Cells.Item(132, 14).Value = -15087.7145  # N132: was -16458.8

# Sheet BSM row 35
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 0  # H35: was 35000
$ws.Cells.Item(35, 10).Value = 0  # J35: was 35000
$ws.Cells.Item(35, 12).Value = 0  # L35: was 35000
$ws.Cells.Item(35, 14).ClearContents()  # N35: was -35620

# Sheet BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 531.7917  # H94: was 551.3043
$ws.Cells.Item(94, 9).Value = 402.7  # I94: was 419.5263
$ws.Cells.Item(94, 11).Value = 402.7  # K94: was 419.5263
$ws.Cells.Item(94, 13).Value = 48.30000000000001  # M94: was 31.47370000000001

# Sheet BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 2183.3333  # H107: was 1712.1875
$ws.Cells.Item(107, 9).Value = 1805.1111  # I107: was 1341.6154
$ws.Cells.Item(107, 11).Value = 1805.1111  # K107: was 1341.6154
$ws.Cells.Item(107, 13).Value = 114.8888999999999  # M107: was 578.3846000000001

# Sheet BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 7397.5806  # H134: was 7604.2334
$ws.Cells.Item(134, 9).Value = 7740  # I134: was 7991.615
$ws.Cells.Item(134, 11).Value = 23220  # K134: was 23974.845
$ws.Cells.Item(134, 13).Value = -20685  # M134: was -21439.845

# Sheet CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3319.4736  # H31: was 3467.2222
$ws.Cells.Item(31, 9).Value = 970  # I31: was 1001
$ws.Cells.Item(31, 11).Value = 970  # K31: was 1001
$ws.Cells.Item(31, 13).Value = -675  # M31: was -706

# Sheet CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 3319.4736  # H34: was 3467.2222
$ws.Cells.Item(34, 9).Value = 970  # I34: was 1001
$ws.Cells.Item(34, 11).Value = 970  # K34: was 1001
$ws.Cells.Item(34, 13).Value = -768  # M34: was -799

# Sheet CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 1999  # H105: was 1999.5
$ws.Cells.Item(105, 9).Value = 1998.5  # I105: was 1999
$ws.Cells.Item(105, 11).Value = 1998.5  # K105: was 1999
$ws.Cells.Item(105, 13).Value = -251.5  # M105: was -252

# Sheet CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2728.0588  # H132: was 2835.3125
$ws.Cells.Item(132, 9).Value = 1136.5  # I132: was 1150.3334
$ws.Cells.Item(132, 11).Value = 3409.5  # K132: was 3451.0002
$ws.Cells.Item(132, 13).Value = -879.5  # M132: was -921.0001999999999

# Sheet CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 1917.3043  # H134: was 1811
$ws.Cells.Item(134, 9).Value = 1738.2778  # I134: was 1622.174
$ws.Cells.Item(134, 10).Value = 2561.8  # J134: was 2679.6
$ws.Cells.Item(134, 11).Value = 5214.8334  # K134: was 4866.522
$ws.Cells.Item(134, 12).Value = 7685.400000000001  # L134: was 8038.799999999999
$ws.Cells.Item(134, 13).Value = -2679.8334  # M134: was -2331.522
$ws.Cells.Item(134, 14).Value = -12755.4  # N134: was -13108.8

# Sheet CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 3966.1667  # H39: was 4574.75
$ws.Cells.Item(39, 10).Value = 3966.1667  # J39: was 4574.75
$ws.Cells.Item(39, 12).Value = 11898.5001  # L39: was 13724.25
$ws.Cells.Item(39, 14).Value = -12486.5001  # N39: was -14312.25

# Sheet CUL row 51
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(51, 8).Value = 750  # H51: was 0
$ws.Cells.Item(51, 9).Value = 750  # I51: was 0
$ws.Cells.Item(51, 11).Value = 2250  # K51: was 0
$ws.Cells.Item(51, 13).Value = -1790  # M51: was None

# Sheet CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 26500.75  # H55: was 21300.6
$ws.Cells.Item(55, 10).Value = 5000  # J55: was 2750
$ws.Cells.Item(55, 12).Value = 15000  # L55: was 8250
$ws.Cells.Item(55, 14).Value = -15354  # N55: was -8604

# Sheet CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 39167.105  # H129: was 52754.5
$ws.Cells.Item(129, 9).Value = 895.2  # I129: was 952
$ws.Cells.Item(129, 10).Value = 52835.645  # J129: was 73475.5
$ws.Cells.Item(129, 11).Value = 2685.6  # K129: was 2856
$ws.Cells.Item(129, 12).Value = 158506.935  # L129: was 220426.5
$ws.Cells.Item(129, 13).Value = 2314.4  # M129: was 2144
$ws.Cells.Item(129, 14).Value = -168506.935  # N129: was -230426.5

# Sheet CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 7474115.5  # H131: was 7257496.5
$ws.Cells.Item(131, 10).Value = 12495.541  # J131: was 12123.302
$ws.Cells.Item(131, 12).Value = 37486.623  # L131: was 36369.906
$ws.Cells.Item(131, 14).Value = -47566.623  # N131: was -46449.906

# Sheet CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 1980  # H132: was 2000
$ws.Cells.Item(132, 10).Value = 2550  # J132: was 3200
$ws.Cells.Item(132, 12).Value = 22950  # L132: was 28800
$ws.Cells.Item(132, 14).Value = -28010  # N132: was -33860

# Sheet GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2445.3845  # H80: was 2482.2222
$ws.Cells.Item(80, 9).Value = 1983.3334  # I80: was 1850
$ws.Cells.Item(80, 10).Value = 2584  # J80: was 2662.8572
$ws.Cells.Item(80, 11).Value = 1983.3334  # K80: was 1850
$ws.Cells.Item(80, 12).Value = 2584  # L80: was 2662.8572
$ws.Cells.Item(80, 13).Value = -985.3334  # M80: was -852
$ws.Cells.Item(80, 14).Value = -4580  # N80: was -4658.8572

# Sheet GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 2445.3845  # H83: was 2482.2222
$ws.Cells.Item(83, 9).Value = 1983.3334  # I83: was 1850
$ws.Cells.Item(83, 10).Value = 2584  # J83: was 2662.8572
$ws.Cells.Item(83, 11).Value = 9916.666999999999  # K83: was 9250
$ws.Cells.Item(83, 12).Value = 12920  # L83: was 13314.286
$ws.Cells.Item(83, 13).Value = -4924.666999999999  # M83: was -4258
$ws.Cells.Item(83, 14).Value = -22904  # N83: was -23298.286

# Sheet GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 545.95  # H97: was 487.47827
$ws.Cells.Item(97, 9).Value = 399.69232  # I97: was 359.6
$ws.Cells.Item(97, 10).Value = 817.5714  # J97: was 727.25
$ws.Cells.Item(97, 11).Value = 399.69232  # K97: was 359.6
$ws.Cells.Item(97, 12).Value = 817.5714  # L97: was 727.25
$ws.Cells.Item(97, 13).Value = 96.30768  # M97: was 136.4
$ws.Cells.Item(97, 14).Value = -1809.5714  # N97: was -1719.25

# Sheet GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2432.5557  # H102: was 2450.8
$ws.Cells.Item(102, 9).Value = 3087.4285  # I102: was 2605.3333
$ws.Cells.Item(102, 10).Value = 2015.8182  # J102: was 2219
$ws.Cells.Item(102, 11).Value = 3087.4285  # K102: was 2605.3333
$ws.Cells.Item(102, 12).Value = 2015.8182  # L102: was 2219
$ws.Cells.Item(102, 13).Value = -1465.4285  # M102: was -983.3332999999998
$ws.Cells.Item(102, 14).Value = -5259.8182  # N102: was -5463

# Sheet GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 4146.409  # H132: was 4290.45
$ws.Cells.Item(132, 9).Value = 3004.8572  # I132: was 3054.6667
$ws.Cells.Item(132, 11).Value = 9014.571599999999  # K132: was 9164.000100000001
$ws.Cells.Item(132, 13).Value = -6484.571599999999  # M132: was -6634.000100000001

# Sheet LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 1060.625  # H93: was 1070.625
$ws.Cells.Item(93, 9).Value = 927.5714  # I93: was 939
$ws.Cells.Item(93, 11).Value = 927.5714  # K93: was 939
$ws.Cells.Item(93, 13).Value = 320.4286  # M93: was 309

# Sheet LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 2341.2856  # H132: was 2537
$ws.Cells.Item(132, 9).Value = 1328.1  # I132: was 1338.8572
$ws.Cells.Item(132, 10).Value = 4874.25  # J132: was 5332.6665
$ws.Cells.Item(132, 11).Value = 3984.3  # K132: was 4016.5716
$ws.Cells.Item(132, 12).Value = 14622.75  # L132: was 15997.9995
$ws.Cells.Item(132, 13).Value = -1454.3  # M132: was -1486.5716
$ws.Cells.Item(132, 14).Value = -19682.75  # N132: was -21057.9995

# Sheet LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 3440.8438  # H136: was 3364.742
$ws.Cells.Item(136, 9).Value = 2252.8572  # I136: was 2300.4546
$ws.Cells.Item(136, 10).Value = 5708.8184  # J136: was 5966.3335
$ws.Cells.Item(136, 11).Value = 6758.571599999999  # K136: was 6901.3638
$ws.Cells.Item(136, 12).Value = 17126.4552  # L136: was 17899.0005
$ws.Cells.Item(136, 13).Value = -4208.571599999999  # M136: was -4351.3638
$ws.Cells.Item(136, 14).Value = -22226.4552  # N136: was -22999.0005

# Sheet LTW row 139
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(139, 8).Value = 45900  # H139: was 46000
$ws.Cells.Item(139, 10).Value = 45900  # J139: was 46000
$ws.Cells.Item(139, 12).Value = 45900  # L139: was 46000
$ws.Cells.Item(139, 14).Value = -56180  # N139: was -56280

# Sheet WVR row 42
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(42, 8).Value = 0  # H42: was 70049
$ws.Cells.Item(42, 10).Value = 0  # J42: was 70049
$ws.Cells.Item(42, 12).Value = 0  # L42: was 70049
$ws.Cells.Item(42, 14).ClearContents()  # N42: was -70805

# Sheet WVR row 48
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(48, 8).Value = 0  # H48: was 44999.5
$ws.Cells.Item(48, 10).Value = 0  # J48: was 44999.5
$ws.Cells.Item(48, 12).Value = 0  # L48: was 44999.5
$ws.Cells.Item(48, 14).ClearContents()  # N48: was -46137.5

# Sheet WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 2500  # H62: was 0
$ws.Cells.Item(62, 9).Value = 3000  # I62: was 0
$ws.Cells.Item(62, 10).Value = 2000  # J62: was 0
$ws.Cells.Item(62, 11).Value = 3000  # K62: was 0
$ws.Cells.Item(62, 12).Value = 2000  # L62: was 0
$ws.Cells.Item(62, 13).Value = -2376  # M62: was None
$ws.Cells.Item(62, 14).Value = -3248  # N62: was None

# Sheet WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(65, 8).Value = 2500  # H65: was 0
$ws.Cells.Item(65, 9).Value = 3000  # I65: was 0
$ws.Cells.Item(65, 10).Value = 2000  # J65: was 0
$ws.Cells.Item(65, 11).Value = 15000  # K65: was 0
$ws.Cells.Item(65, 12).Value = 10000  # L65: was 0
$ws.Cells.Item(65, 13).Value = -11880  # M65: was None
$ws.Cells.Item(65, 14).Value = -16240  # N65: was None

# Sheet WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1989.0416  # H136: was 2418.7805
$ws.Cells.Item(136, 9).Value = 1685.359  # I136: was 1954.875
$ws.Cells.Item(136, 10).Value = 3305  # J136: was 4068.2222
$ws.Cells.Item(136, 11).Value = 5056.076999999999  # K136: was 5864.625
$ws.Cells.Item(136, 12).Value = 9915  # L136: was 12204.6666
$ws.Cells.Item(136, 13).Value = -2506.076999999999  # M136: was -3314.625
$ws.Cells.Item(136, 14).Value = -15015  # N136: was -17304.6666

# Sheet WVR row 139
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(139, 8).Value = 59942.855  # H139: was 56353.824
$ws.Cells.Item(139, 10).Value = 59942.855  # J139: was 56353.824
$ws.Cells.Item(139, 12).Value = 59942.855  # L139: was 56353.824
$ws.Cells.Item(139, 14).Value = -70222.85500000001  # N139: was -66633.82399999999

Write-Host "All updates applied"